# Update countries & provincias Spain
# Applies the diff: two country label re-sorts (Suiza/Paises Bajos and
# Moldavia/Croacia swap position based on updated case counts), refreshed
# case numbers for several rows, and a new "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 / 15: Paises Bajos overtakes Suiza -------------------------
$ws.Cells.Item(14, 1).Value = "Paises Bajos"
$ws.Cells.Item(14, 2).Value = 25587
$ws.Cells.Item(14, 3).Value = 1174
$ws.Cells.Item(14, 4).Value = 250
$ws.Cells.Item(14, 5).Value = 22600
$ws.Cells.Item(14, 6).Value = 1384
$ws.Cells.Item(14, 7).Value = 94
$ws.Cells.Item(14, 8).Value = 2737

$ws.Cells.Item(15, 1).Value = "Suiza"
$ws.Cells.Item(15, 2).Value = 25300
$ws.Cells.Item(15, 3).Value = 193
$ws.Cells.Item(15, 4).Value = 12100
$ws.Cells.Item(15, 5).Value = 12164
$ws.Cells.Item(15, 6).Value = 386
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 1036

# --- Row 18: Portugal, updated counts -----------------------------------
$ws.Cells.Item(18, 2).Value = 16585
$ws.Cells.Item(18, 3).Value = 598
$ws.Cells.Item(18, 4).Value = 277
$ws.Cells.Item(18, 5).Value = 15804
$ws.Cells.Item(18, 6).Value = 228
$ws.Cells.Item(18, 7).Value = 34
$ws.Cells.Item(18, 8).Value = 504

# --- Row 34: Dinamarca, updated counts ----------------------------------
$ws.Cells.Item(34, 2).Value = 6174
$ws.Cells.Item(34, 3).Value = 178
$ws.Cells.Item(34, 4).Value = 2123
$ws.Cells.Item(34, 5).Value = 3778
$ws.Cells.Item(34, 6).Value = 104
$ws.Cells.Item(34, 7).Value = 13
$ws.Cells.Item(34, 8).Value = 273

# --- Row 54: Argentina, updated counts ----------------------------------
$ws.Cells.Item(54, 5).Value = 1612
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 90

# --- Row 61 / 62: Croacia overtakes Moldavia -----------------------------
$ws.Cells.Item(61, 1).Value = "Croacia"
$ws.Cells.Item(61, 2).Value = 1600
$ws.Cells.Item(61, 3).Value = 66
$ws.Cells.Item(61, 4).Value = 373
$ws.Cells.Item(61, 5).Value = 1204
$ws.Cells.Item(61, 6).Value = 34
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 23

$ws.Cells.Item(62, 1).Value = "Moldavia"
$ws.Cells.Item(62, 2).Value = 1560
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 75
$ws.Cells.Item(62, 5).Value = 1455
$ws.Cells.Item(62, 6).Value = 80
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 30

# --- Row 111: Vietnam, updated counts -----------------------------------
$ws.Cells.Item(111, 2).Value = 260
$ws.Cells.Item(111, 3).Value = 2
$ws.Cells.Item(111, 5).Value = 116

# --- Row 159: Haiti, updated counts --------------------------------------
$ws.Cells.Item(159, 5).Value = 30
$ws.Cells.Item(159, 7).Value = 1
$ws.Cells.Item(159, 8).Value = 3

# --- Header timestamp -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Abril de 2020 a las 14:22"
